# Natmi following Dr Hou advice
# Add "ECs" as an additional Sending-cluster group (alongside the existing
# FAPs / sCs clusters) for the Lama1 -> Itgb8 ligand-receptor pair, so the
# sheet now covers the full 3x3 cross of {ECs, FAPs, sCs} sending clusters
# against {ECs, FAPs, sCs} target clusters (rows 2-10) instead of the
# previous 2x3 cross (rows 2-7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Lama1"
$ws.Cells.Item(2,3).Value = "Itgb8"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.02283333333333333
$ws.Cells.Item(2,8).Value = 0.06850000000000001
$ws.Cells.Item(2,9).Value = 0.05477178157813095
$ws.Cells.Item(2,10).Value = 0.05477178157813096
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.3116673333333334
$ws.Cells.Item(2,14).Value = 0.935002
$ws.Cells.Item(2,15).Value = 0.0414413620607491
$ws.Cells.Item(2,16).Value = 0.0414413620607491
$ws.Cells.Item(2,17).Value = 0.007116404111111111
$ws.Cells.Item(2,18).Value = 0.064047637
$ws.Cells.Item(2,19).Value = 0.002269817231091592
$ws.Cells.Item(2,20).Value = 0.002269817231091592
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Lama1"
$ws.Cells.Item(3,3).Value = "Itgb8"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.02283333333333333
$ws.Cells.Item(3,8).Value = 0.06850000000000001
$ws.Cells.Item(3,9).Value = 0.05477178157813095
$ws.Cells.Item(3,10).Value = 0.05477178157813096
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 3.794584666666667
$ws.Cells.Item(3,14).Value = 11.383754
$ws.Cells.Item(3,15).Value = 0.5045532214096876
$ws.Cells.Item(3,16).Value = 0.5045532214096876
$ws.Cells.Item(3,17).Value = 0.08664301655555556
$ws.Cells.Item(3,18).Value = 0.779787149
$ws.Cells.Item(3,19).Value = 0.02763527883759375
$ws.Cells.Item(3,20).Value = 0.02763527883759376
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Lama1"
$ws.Cells.Item(4,3).Value = "Itgb8"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.02283333333333333
$ws.Cells.Item(4,8).Value = 0.06850000000000001
$ws.Cells.Item(4,9).Value = 0.05477178157813095
$ws.Cells.Item(4,10).Value = 0.05477178157813096
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 3.414430666666667
$ws.Cells.Item(4,14).Value = 10.243292
$ws.Cells.Item(4,15).Value = 0.4540054165295632
$ws.Cells.Item(4,16).Value = 0.4540054165295633
$ws.Cells.Item(4,17).Value = 0.07796283355555556
$ws.Cells.Item(4,18).Value = 0.7016655020000001
$ws.Cells.Item(4,19).Value = 0.0248666855094456
$ws.Cells.Item(4,20).Value = 0.02486668550944561
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Lama1"
$ws.Cells.Item(5,3).Value = "Itgb8"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.3395593333333333
$ws.Cells.Item(5,8).Value = 1.018678
$ws.Cells.Item(5,9).Value = 0.8145227578751427
$ws.Cells.Item(5,10).Value = 0.8145227578751428
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.3116673333333334
$ws.Cells.Item(5,14).Value = 0.935002
$ws.Cells.Item(5,15).Value = 0.0414413620607491
$ws.Cells.Item(5,16).Value = 0.0414413620607491
$ws.Cells.Item(5,17).Value = 0.1058295519284445
$ws.Cells.Item(5,18).Value = 0.952465967356
$ws.Cells.Item(5,19).Value = 0.03375493251582366
$ws.Cells.Item(5,20).Value = 0.03375493251582366
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Lama1"
$ws.Cells.Item(6,3).Value = "Itgb8"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0.3395593333333333
$ws.Cells.Item(6,8).Value = 1.018678
$ws.Cells.Item(6,9).Value = 0.8145227578751427
$ws.Cells.Item(6,10).Value = 0.8145227578751428
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 3.794584666666667
$ws.Cells.Item(6,14).Value = 11.383754
$ws.Cells.Item(6,15).Value = 0.5045532214096876
$ws.Cells.Item(6,16).Value = 0.5045532214096876
$ws.Cells.Item(6,17).Value = 1.288486639690222
$ws.Cells.Item(6,18).Value = 11.596379757212
$ws.Cells.Item(6,19).Value = 0.4109700813974063
$ws.Cells.Item(6,20).Value = 0.4109700813974063
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Lama1"
$ws.Cells.Item(7,3).Value = "Itgb8"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 0.3395593333333333
$ws.Cells.Item(7,8).Value = 1.018678
$ws.Cells.Item(7,9).Value = 0.8145227578751427
$ws.Cells.Item(7,10).Value = 0.8145227578751428
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 3.414430666666667
$ws.Cells.Item(7,14).Value = 10.243292
$ws.Cells.Item(7,15).Value = 0.4540054165295632
$ws.Cells.Item(7,16).Value = 0.4540054165295633
$ws.Cells.Item(7,17).Value = 1.159401800886222
$ws.Cells.Item(7,18).Value = 10.434616207976
$ws.Cells.Item(7,19).Value = 0.3697977439619127
$ws.Cells.Item(7,20).Value = 0.3697977439619128
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Lama1"
$ws.Cells.Item(8,3).Value = "Itgb8"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.05448866666666666
$ws.Cells.Item(8,8).Value = 0.163466
$ws.Cells.Item(8,9).Value = 0.1307054605467263
$ws.Cells.Item(8,10).Value = 0.1307054605467264
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 0.3116673333333334
$ws.Cells.Item(8,14).Value = 0.935002
$ws.Cells.Item(8,15).Value = 0.0414413620607491
$ws.Cells.Item(8,16).Value = 0.0414413620607491
$ws.Cells.Item(8,17).Value = 0.01698233743688889
$ws.Cells.Item(8,18).Value = 0.152841036932
$ws.Cells.Item(8,19).Value = 0.005416612313833842
$ws.Cells.Item(8,20).Value = 0.005416612313833843
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Lama1"
$ws.Cells.Item(9,3).Value = "Itgb8"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.05448866666666666
$ws.Cells.Item(9,8).Value = 0.163466
$ws.Cells.Item(9,9).Value = 0.1307054605467263
$ws.Cells.Item(9,10).Value = 0.1307054605467264
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 3.794584666666667
$ws.Cells.Item(9,14).Value = 11.383754
$ws.Cells.Item(9,15).Value = 0.5045532214096876
$ws.Cells.Item(9,16).Value = 0.5045532214096876
$ws.Cells.Item(9,17).Value = 0.2067618590404444
$ws.Cells.Item(9,18).Value = 1.860856731364
$ws.Cells.Item(9,19).Value = 0.06594786117468759
$ws.Cells.Item(9,20).Value = 0.06594786117468761
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Lama1"
$ws.Cells.Item(10,3).Value = "Itgb8"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.05448866666666666
$ws.Cells.Item(10,8).Value = 0.163466
$ws.Cells.Item(10,9).Value = 0.1307054605467263
$ws.Cells.Item(10,10).Value = 0.1307054605467264
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 3.414430666666667
$ws.Cells.Item(10,14).Value = 10.243292
$ws.Cells.Item(10,15).Value = 0.4540054165295632
$ws.Cells.Item(10,16).Value = 0.4540054165295633
$ws.Cells.Item(10,17).Value = 0.1860477744524444
$ws.Cells.Item(10,18).Value = 1.674429970072
$ws.Cells.Item(10,19).Value = 0.05934098705820488
$ws.Cells.Item(10,20).Value = 0.0593409870582049
